# Applies a row-permutation of the date/quality/volume/price columns
# (D, I, J, K, L, M, P) across rows 2..43 of the active worksheet.
# The mapping below says: new row <key> receives the values that
# currently (before the edit) live in row <value>.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 6
    3  = 20
    4  = 31
    5  = 41
    6  = 29
    7  = 16
    8  = 43
    9  = 27
    10 = 14
    11 = 2
    12 = 42
    13 = 19
    14 = 18
    15 = 28
    16 = 15
    17 = 3
    18 = 13
    19 = 11
    20 = 12
    21 = 36
    22 = 26
    23 = 32
    24 = 34
    25 = 17
    26 = 10
    27 = 23
    28 = 7
    29 = 35
    30 = 9
    31 = 22
    32 = 37
    33 = 4
    34 = 25
    35 = 5
    36 = 38
    37 = 39
    38 = 40
    39 = 21
    40 = 8
    41 = 30
    42 = 33
    43 = 24
}

$cols = @(4, 9, 10, 11, 12, 13, 16)   # D, I, J, K, L, M, P

# Snapshot every source cell's value first, since sources and
# destinations overlap (it's a permutation, not a copy from a
# separate range).
$snapshot = @{}
foreach ($row in 2..43) {
    foreach ($col in $cols) {
        $snapshot["$row`_$col"] = $ws.Cells.Item($row, $col).Value()
    }
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    foreach ($col in $cols) {
        $val = $snapshot["$oldRow`_$col"]
        $ws.Cells.Item($newRow, $col).Value = $val
    }
}
